{"js": "// Merge the split \"<id>...</id>\" runs (opening tag / value / closing tag,\n// each with their own run formatting) back into a single run per\n// occurrence, e.g. \"<id>p111r_1</id>\" for p111r_1 .. p111r_5.\n// The standalone \"<id>fig_p111r_1</id>\" run triple is intentionally left\n// untouched, matching the source diff.\n\nconst body = context.document.body;\n\nfor (let n = 1; n <= 5; n++) {\n  const needle = `<id>p111r_${n}</id>`;\n\n  const results = body.search(needle, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    // Replacing the whole matched range (which spans the 3 original runs)\n    // with its own text collapses it into a single run using the first\n    // run's formatting (Courier New / 7f6000 / 18pt) - exactly matching\n    // the target edit.\n    results.items[0].insertText(needle, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Merge the split \"<id>...</id>\" runs (tag/value/tag) back into a single\n# run per occurrence, e.g. \"<id>p111r_1</id>\" for p111r_1 .. p111r_5.\n# (The standalone \"<id>fig_p111r_1</id>\" paragraph is left untouched, as\n# in the source diff.)\n\n$d = $word.ActiveDocument\n\nfor ($n = 1; $n -le 5; $n++) {\n    $needle = \"<id>p111r_$n</id>\"\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $needle\n    $find.Forward = $true\n    $found = $find.Execute()\n\n    if ($found) {\n        $rng = $find.Parent\n\n        # Skip past the literal \"<id>\" (4 chars) so the remaining range -\n        # the id value plus the closing \"</id>\" tag - is deleted and\n        # reinserted immediately after the still-intact opening \"<id>\"\n        # run. The new text then naturally inherits that run's\n        # formatting (Courier New / 7f6000 / 18pt), merging all three\n        # original runs into one, exactly like the target edit.\n        $rng.MoveStart(1, 4)\n        $tail = $rng.Text\n        $rng.Delete()\n        $rng.InsertAfter($tail)\n    }\n}\n"}
